$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stock Report")

# Reference cell carrying the date-style (s="4") used for columns M and U
$dateStyleRef = "U6"

# --- Row 7: replace the single stringified-array cell in B7 with the
#     real per-column values it represented ---
$ws.Range("B7:V7").ClearContents()
$ws.Range("B7:V7").Style = "Normal"
$ws.Range($dateStyleRef).Copy()
$ws.Range("M7").PasteSpecial(-4122)
$ws.Range($dateStyleRef).Copy()
$ws.Range("U7").PasteSpecial(-4122)
$ws.Range("W7").Value = "FLOORS-(F)"
$ws.Range("X7").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y7").Value = "FLOOR BOARD DIRTY BY DUST"
$ws.Range("Z7").Style = "Normal"
$ws.Range("AA7").Style = "Normal"
$ws.Range("AB7").Style = "Normal"

# --- Row 9: replace the single stringified-array cell in B9 with the
#     real per-column values it represented ---
$ws.Range("B9:V9").ClearContents()
$ws.Range("B9:V9").Style = "Normal"
$ws.Range($dateStyleRef).Copy()
$ws.Range("M9").PasteSpecial(-4122)
$ws.Range($dateStyleRef).Copy()
$ws.Range("U9").PasteSpecial(-4122)
$ws.Range("W9").Value = "FLOORS-(F)"
$ws.Range("X9").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y9").Value = "F/BOARD DIRTY BY DUST."
$ws.Range("Z9").Style = "Normal"
$ws.Range("AA9").Style = "Normal"
$ws.Range("AB9").Style = "Normal"

# --- Row 10: replace the single stringified-array cell in B10 with the
#     real per-column values it represented ---
$ws.Range("B10:V10").ClearContents()
$ws.Range("B10:V10").Style = "Normal"
$ws.Range($dateStyleRef).Copy()
$ws.Range("M10").PasteSpecial(-4122)
$ws.Range($dateStyleRef).Copy()
$ws.Range("U10").PasteSpecial(-4122)
$ws.Range("W10").Value = "DOORS-(D)"
$ws.Range("X10").Value = "Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)"
$ws.Range("Y10").Value = "R/DOOR BTM GASKET CUT 06`""
$ws.Range("Z10").Style = "Normal"
$ws.Range("AA10").Style = "Normal"
$ws.Range("AB10").Style = "Normal"

# --- Row 12: replace the single stringified-array cell in B12 with the
#     real per-column values it represented ---
$ws.Range("B12:V12").ClearContents()
$ws.Range("B12:V12").Style = "Normal"
$ws.Range($dateStyleRef).Copy()
$ws.Range("M12").PasteSpecial(-4122)
$ws.Range($dateStyleRef).Copy()
$ws.Range("U12").PasteSpecial(-4122)
$ws.Range("W12").Value = "FLOORS-(F)"
$ws.Range("X12").Value = "FLOOR BOARD-(FLOOR BOARD)"
$ws.Range("Y12").Value = "F/B DIRTY BY DUST ."
$ws.Range("Z12").Style = "Normal"
$ws.Range("AA12").Style = "Normal"
$ws.Range("AB12").Style = "Normal"

# --- Row 14: replace the single stringified-array cell in B14 with the
#     real per-column values it represented ---
$ws.Range("B14:V14").ClearContents()
$ws.Range("B14:V14").Style = "Normal"
$ws.Range($dateStyleRef).Copy()
$ws.Range("M14").PasteSpecial(-4122)
$ws.Range($dateStyleRef).Copy()
$ws.Range("U14").PasteSpecial(-4122)
$ws.Range("W14").Value = "FLOORS-(F)"
$ws.Range("X14").Value = "FLOOR BOARD-(FLOOR BOARD)"
$ws.Range("Y14").Value = "FLOOR BOARD DIRTY BY DUST"
$ws.Range("Z14").Style = "Normal"
$ws.Range("AA14").Style = "Normal"
$ws.Range("AB14").Style = "Normal"

# --- Row 16: replace the single stringified-array cell in B16 with the
#     real per-column values it represented ---
$ws.Range("B16:V16").ClearContents()
$ws.Range("B16:V16").Style = "Normal"
$ws.Range($dateStyleRef).Copy()
$ws.Range("M16").PasteSpecial(-4122)
$ws.Range($dateStyleRef).Copy()
$ws.Range("U16").PasteSpecial(-4122)
$ws.Range("W16").Value = "FLOORS-(F)"
$ws.Range("X16").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y16").Value = "RF FLOOR BOARD DIRTY BY DUST & ODOUR ."
$ws.Range("Z16").Style = "Normal"
$ws.Range("AA16").Style = "Normal"
$ws.Range("AB16").Style = "Normal"

# --- Row 18: replace the single stringified-array cell in B18 with the
#     real per-column values it represented ---
$ws.Range("B18:V18").ClearContents()
$ws.Range("B18:V18").Style = "Normal"
$ws.Range($dateStyleRef).Copy()
$ws.Range("M18").PasteSpecial(-4122)
$ws.Range($dateStyleRef).Copy()
$ws.Range("U18").PasteSpecial(-4122)
$ws.Range("W18").Value = "FLOORS-(F)"
$ws.Range("X18").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y18").Value = "RF FLOOR BOARD DIRTY BY DUST & ODOUR ."
$ws.Range("Z18").Style = "Normal"
$ws.Range("AA18").Style = "Normal"
$ws.Range("AB18").Style = "Normal"

# --- Row 20: replace the single stringified-array cell in B20 with the
#     real per-column values it represented ---
$ws.Range("B20:V20").ClearContents()
$ws.Range("B20:V20").Style = "Normal"
$ws.Range($dateStyleRef).Copy()
$ws.Range("M20").PasteSpecial(-4122)
$ws.Range($dateStyleRef).Copy()
$ws.Range("U20").PasteSpecial(-4122)
$ws.Range("W20").Value = "FLOORS-(F)"
$ws.Range("X20").Value = "Threshold plate-(Threshold plate)"
$ws.Range("Y20").Value = "RF FLOOR BOARD DIRTY BY DUST & ODOUR ."
$ws.Range("Z20").Style = "Normal"
$ws.Range("AA20").Style = "Normal"
$ws.Range("AB20").Style = "Normal"

$excel.CutCopyMode = $false

# --- Column width / visibility adjustments ---
# Column B no longer needs to hold the giant stringified array, so its
# bestFit width shrinks back down (still hidden).
$ws.Columns.Item(2).ColumnWidth = 7.857142857142857
$ws.Columns.Item(2).Hidden = $true
# Columns X and Y now contain the longest strings in the sheet, so their
# bestFit widths grow to fit.
$ws.Columns.Item(24).ColumnWidth = 46.42857142857143
$ws.Columns.Item(25).ColumnWidth = 42.0
